# Add Trade #4 (closed/open info as of 2026-02-16 22:52:20) to both the
# "All Trades" summary sheet and the per-strategy "base_strategy" sheet.
# Both sheets hold identical trade logs, so the new row is appended to each.

$wb = $excel.ActiveWorkbook

$sheetNames = @("All Trades", "base_strategy")

foreach ($sheetName in $sheetNames) {
    $ws = $wb.Worksheets.Item($sheetName)

    $row = 5

    $ws.Cells.Item($row, 1).Value = 4                       # Trade #
    # Leading apostrophe forces literal text so Excel doesn't reinterpret
    # the yyyy-mm-dd string as a date serial (matches how the existing
    # Date column cells above it are stored as plain text).
    $ws.Cells.Item($row, 2).Value = "'2026-02-16"             # Date
    $ws.Cells.Item($row, 3).Value = "22:52:20"                # Time
    $ws.Cells.Item($row, 4).Value = "base_strategy"           # Strategy
    $ws.Cells.Item($row, 5).Value = "DOWN"                     # Side
    $ws.Cells.Item($row, 6).Value = 49.999998                  # Entry Price
    $ws.Cells.Item($row, 7).Value = ""                          # Exit Price (blank - trade still open)
    $ws.Cells.Item($row, 8).Value = "OPEN"                       # Status
    $ws.Cells.Item($row, 9).Value = 0                            # P&L %
    $ws.Cells.Item($row, 10).Value = 0                           # P&L $
    $ws.Cells.Item($row, 11).Value = 100                         # Capital After
    $ws.Cells.Item($row, 12).Value = 0                           # Entry Slippage (bps)
    $ws.Cells.Item($row, 13).Value = 0                           # Exit Slippage (bps)
    $ws.Cells.Item($row, 14).Value = 0.6                         # Confidence
    $ws.Cells.Item($row, 15).Value = "Normal spread capture: 19600 bps"  # Entry Reason
    $ws.Cells.Item($row, 16).Value = ""                           # Exit Reason (blank - trade still open)
    $ws.Cells.Item($row, 17).Value = 0                            # Duration (min)
}
